# Commit 3: Final Drafts of Math and CS Visualizations in both BGWP and Color
# Added Math Dep to dataset; Added Color to CS Visualization; Completed Math Visualization

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Math Department")

# --- Add new header labels to the Math Department sheet ---
# Order matters so that the shared-strings table gets the same index
# assignment as the target workbook (Math, Applied Math, Math/Econ,
# Secondary Education, Applied Sciences, Probability and Statistics).
$ws2.Range("A1").Value = "Quarter"
$ws2.Range("E1").Value = "Math"
$ws2.Range("C1").Value = "Applied Math"
$ws2.Range("F1").Value = "Math/Econ"
$ws2.Range("D1").Value = "Secondary Education"
$ws2.Range("G1").Value = "Applied Sciences"
$ws2.Range("H1").Value = "Probability and Statistics"

# --- Fill in the new numeric data (columns C:H, rows 2:15) ---
$data = @(
    @(73, 39, 166, 117, 21, 40),
    @(80, 36, 152, 118, 17, 41),
    @(86, 37, 131, 114, 14, 40),
    @(125, 40, 145, 136, 21, 37),
    @(135, 36, 124, 147, 18, 40),
    @(136, 35, 110, 145, 17, 43),
    @(185, 35, 149, 167, 31, 53),
    @(200, 34, 130, 164, 26, 60),
    @(189, 37, 121, 163, 24, 62),
    @(299, 42, 182, 156, 44, 96),
    @(313, 36, 155, 158, 45, 111),
    @(311, 31, 147, 170, 31, 119),
    @(417, 40, 219, 192, 46, 162),
    @(409, 37, 211, 241, 46, 160)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws2.Cells.Item($row, 3).Value = $vals[0]
    $ws2.Cells.Item($row, 4).Value = $vals[1]
    $ws2.Cells.Item($row, 5).Value = $vals[2]
    $ws2.Cells.Item($row, 6).Value = $vals[3]
    $ws2.Cells.Item($row, 7).Value = $vals[4]
    $ws2.Cells.Item($row, 8).Value = $vals[5]
}

# --- Column E width adjustment (target stored width 8.7265625 chars) ---
# The host's ColumnWidth setter quantizes the resulting stored width to
# multiples of 1/6 character (5/6 fixed padding + nearest-1/6 of the
# requested value), so 47/6 = 7.8333... is the input that lands on the
# closest achievable stored width (52/6 = 8.666...) to the target.
$ws2.Range("E1").ColumnWidth = 7.8333333333

# --- Activate the Math Department sheet (moves tabSelected / activeTab) ---
$ws2.Activate()

# --- Selection: whole column G is selected on the Math Department sheet ---
$ws2.Range("G:G").Select() | Out-Null
